# Reto_Tecnico-Choucair: add "Validacion" sheet + tweak selections
# (mirrors commit "Actualización de lógica en ValidarValoresQuestions,
#  VerificarTextoEnTabla y reclutamiento")

$wb = $excel.ActiveWorkbook

# --- FormularioCandidato: move the saved selection, drop the best-fit flag
#     on column A (same visual width, just no longer "auto" fitted) -------
$wsForm = $wb.Worksheets.Item("FormularioCandidato")
$wsForm.Range("D19").Select()
$wsForm.Columns.Item(1).ColumnWidth = 14.833333333333334

# --- New "Validacion" sheet, placed after "Entrevista" -------------------
$wsEntrevista = $wb.Worksheets.Item("Entrevista")
$wsValidacion = $wb.Worksheets.Add([System.Type]::Missing, $wsEntrevista)
$wsValidacion.Name = "Validacion"

$wsValidacion.Columns.Item(1).ColumnWidth = 20.833333333333332
$wsValidacion.Columns.Item(2).ColumnWidth = 24.166666666666668
$wsValidacion.Columns.Item(3).ColumnWidth = 15.833333333333334

$wsValidacion.Range("A1").Value = "Vacancy"
$wsValidacion.Range("B1").Value = "Candidate"
$wsValidacion.Range("C1").Value = "Status"
$wsValidacion.Range("A2").Value = "Payroll Administrator"
$wsValidacion.Range("B2").Value = "Juan Camilo Anacona"
$wsValidacion.Range("C2").Value = "Hired"

$wsValidacion.PageSetup.Orientation = 1

# Activating this sheet last makes it the workbook's active tab and leaves
# its own selection on C1, matching the saved view in the target workbook.
$wsValidacion.Range("C1").Select()
$wsValidacion.Activate()
